{"js": "// Replace each three-digit x one-digit multiplication \"problem=answer\" string\n// with its new value. Every \"old\" value is unique in the document, so a\n// search-and-replace-first-hit per pair is unambiguous. The one pair whose\n// \"new\" text equals another pair's \"old\" text (190x6=1140) is ordered first,\n// so it is consumed before it is (re)created later in the list.\nconst replacements = [\n  [\"190\u00d76=1140\", \"884\u00d72=1768\"],\n  [\"579\u00d78=4632\", \"960\u00d73=2880\"],\n  [\"502\u00d73=1506\", \"405\u00d78=3240\"],\n  [\"277\u00d77=1939\", \"414\u00d77=2898\"],\n  [\"455\u00d77=3185\", \"190\u00d76=1140\"],\n  [\"222\u00d75=1110\", \"701\u00d74=2804\"],\n  [\"779\u00d74=3116\", \"371\u00d77=2597\"],\n  [\"101\u00d73=303\", \"943\u00d77=6601\"],\n  [\"236\u00d77=1652\", \"507\u00d73=1521\"],\n  [\"398\u00d78=3184\", \"222\u00d78=1776\"],\n  [\"124\u00d79=1116\", \"975\u00d75=4875\"],\n  [\"502\u00d79=4518\", \"284\u00d78=2272\"],\n  [\"827\u00d73=2481\", \"539\u00d76=3234\"],\n  [\"693\u00d73=2079\", \"852\u00d74=3408\"],\n  [\"921\u00d75=4605\", \"817\u00d77=5719\"],\n  [\"387\u00d73=1161\", \"171\u00d77=1197\"],\n  [\"457\u00d76=2742\", \"866\u00d79=7794\"],\n  [\"790\u00d74=3160\", \"714\u00d75=3570\"],\n  [\"103\u00d72=206\", \"765\u00d79=6885\"],\n  [\"879\u00d74=3516\", \"410\u00d78=3280\"],\n  [\"258\u00d79=2322\", \"940\u00d73=2820\"],\n  [\"841\u00d72=1682\", \"605\u00d75=3025\"],\n  [\"371\u00d75=1855\", \"424\u00d75=2120\"],\n  [\"840\u00d77=5880\", \"822\u00d76=4932\"],\n  [\"604\u00d73=1812\", \"895\u00d74=3580\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Each pair replaces the unique \"problem=answer\" text of one table cell.\n# \"190x6=1140\" is both an old value (last data row) and a new value\n# (row 1, col 4), so its consuming replacement is ordered first to avoid\n# touching the wrong cell once the new \"190x6=1140\" text exists.\n$pairs = @(\n    @{Old = \"190\u00d76=1140\"; New = \"884\u00d72=1768\"},\n    @{Old = \"579\u00d78=4632\"; New = \"960\u00d73=2880\"},\n    @{Old = \"502\u00d73=1506\"; New = \"405\u00d78=3240\"},\n    @{Old = \"277\u00d77=1939\"; New = \"414\u00d77=2898\"},\n    @{Old = \"455\u00d77=3185\"; New = \"190\u00d76=1140\"},\n    @{Old = \"222\u00d75=1110\"; New = \"701\u00d74=2804\"},\n    @{Old = \"779\u00d74=3116\"; New = \"371\u00d77=2597\"},\n    @{Old = \"101\u00d73=303\"; New = \"943\u00d77=6601\"},\n    @{Old = \"236\u00d77=1652\"; New = \"507\u00d73=1521\"},\n    @{Old = \"398\u00d78=3184\"; New = \"222\u00d78=1776\"},\n    @{Old = \"124\u00d79=1116\"; New = \"975\u00d75=4875\"},\n    @{Old = \"502\u00d79=4518\"; New = \"284\u00d78=2272\"},\n    @{Old = \"827\u00d73=2481\"; New = \"539\u00d76=3234\"},\n    @{Old = \"693\u00d73=2079\"; New = \"852\u00d74=3408\"},\n    @{Old = \"921\u00d75=4605\"; New = \"817\u00d77=5719\"},\n    @{Old = \"387\u00d73=1161\"; New = \"171\u00d77=1197\"},\n    @{Old = \"457\u00d76=2742\"; New = \"866\u00d79=7794\"},\n    @{Old = \"790\u00d74=3160\"; New = \"714\u00d75=3570\"},\n    @{Old = \"103\u00d72=206\"; New = \"765\u00d79=6885\"},\n    @{Old = \"879\u00d74=3516\"; New = \"410\u00d78=3280\"},\n    @{Old = \"258\u00d79=2322\"; New = \"940\u00d73=2820\"},\n    @{Old = \"841\u00d72=1682\"; New = \"605\u00d75=3025\"},\n    @{Old = \"371\u00d75=1855\"; New = \"424\u00d75=2120\"},\n    @{Old = \"840\u00d77=5880\"; New = \"822\u00d76=4932\"},\n    @{Old = \"604\u00d73=1812\"; New = \"895\u00d74=3580\"}\n)\n\nforeach ($p in $pairs) {\n  $range = $d.Content\n  $found = $range.Find.Execute($p.Old, $true, $false, $false, $false, $false, $true, 1, $false, $p.New, 2)\n  if (-not $found) {\n    throw \"Text not found: $($p.Old)\"\n  }\n}\n"}
